$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")
$ws.Activate() | Out-Null

# Header row: AM4 gets the "Julio" month label (new shared string)
$ws.Range("AM4").Value = "Julio"

# Group "Presentar y aprovación de proyecto." (row 8) rolls up rows 9:11
$ws.Range("F8").Formula = "=SUM(F9:F11)"
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = 3

# Group "Planificación" (row 12) rolls up rows 13:17
$ws.Range("F12").Formula = "=SUM(F13:F17)"
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 3
$ws.Range("F16").Value = 3
$ws.Range("F17").Value = 3

# Group "Desarrollo" (row 18) rolls up rows 19:22
$ws.Range("F18").Formula = "=SUM(F19:F22)"
$ws.Range("F19").Value = 3
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 3
$ws.Range("F22").Value = 10

# Group "Pruebas de funcionamiento y Entrega" (row 23) rolls up rows 24:27
$ws.Range("F23").Formula = "=SUM(F24:F27)"
$ws.Range("F24").Value = 3
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 3
$ws.Range("F27").Value = 1

# The Gantt timeline window shifted three columns narrower (AO:AQ now hidden)
$ws.Range("AO1:AQ1").EntireColumn.Hidden = $true

# Final selection left on F28 after entering the last value
$ws.Range("F28").Select() | Out-Null
